$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(9).Insert()
Write-Output $ws.Range("Q10").Style.Name
$ws.Range("Q9").Style = $ws.Range("Q10").Style
